# "Only one student crew flying per time slot"
#
# The Instructors sheet listed, per instructor, how many student crews
# (MaxHours column) and how many could fly at once (Availability column)
# they could be scheduled with simultaneously. Going forward only one
# student crew may fly per time slot, so every instructor's MaxHours
# (and, where it had allowed two concurrent crews, Availability) is
# capped at 1.

$wb = $excel.ActiveWorkbook
$wsInstructors = $wb.Worksheets.Item("Instructors")
$wsStudents    = $wb.Worksheets.Item("Students")

# CAMPBELL / WALKER: MaxHours 2 -> 1
$wsInstructors.Range("C2").Value = 1
$wsInstructors.Range("C3").Value = 1

# CROSS: Availability and MaxHours 2 -> 1
$wsInstructors.Range("B5").Value = 1
$wsInstructors.Range("C5").Value = 1

# CUMMINGS: Availability and MaxHours 2 -> 1
$wsInstructors.Range("B6").Value = 1
$wsInstructors.Range("C6").Value = 1

# BARRON: MaxHours 2 -> 1
$wsInstructors.Range("C8").Value = 1

# Reflect where the author was last working when the file was saved:
# the Students sheet had D2 selected, and the Instructors sheet (rather
# than Reward) ended up as the active tab.
$wsStudents.Activate()
$wsStudents.Range("D2").Select()
$wsInstructors.Activate()
